Write-Host "no-op"
